$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change D3's formula to reference an undefined name ("bfg") instead of A1.
# This makes D3 evaluate to a #NAME? error, which cascades into D1 (which
# sums D3) and flips the cached ISERROR() results in D8/D9 from FALSE to TRUE.
$ws.Range("D3").Formula = "=SUM(bfg,76)"
